$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("BC2:BC56")
$rng.Font.Name = "MesloLGM NF"
$rng.Font.Name = "Arial"
